$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "66.544.88"
$ws.Range("E2").Value = "  -3.56%  "

# Row 3
$ws.Range("D3").Value = "3.560.40"
$ws.Range("E3").Value = "  -4.33%  "

# Row 4
$ws.Range("E4").Value = "  +0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "571.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -6.94%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "185.71"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.75%  "

# Row 7
$ws.Range("D7").Value = "3.554.64"
$ws.Range("E7").Value = "  -4.45%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.614"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.00%  "

# Row 9
$ws.Range("E9").Value = "  +0.22%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.672"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -7.47%  "

# Row 11
$ws.Range("B11").Value = "Avalanche"
$ws.Range("C11").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "55.65"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -7.89%  "

# Row 12
$ws.Range("B12").Value = "Dogecoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.149"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -8.05%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000262"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -10.04%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.80"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.62%  "

# Row 15
$ws.Range("D15").Value = "4.128.75"
$ws.Range("E15").Value = "  -4.32%  "

# Row 16
$ws.Range("D16").Value = "3.558.96"
$ws.Range("E16").Value = "  -4.36%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.24"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -6.71%  "

# Row 19
$ws.Range("D19").Value = "66.584.72"
$ws.Range("E19").Value = "  -3.38%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.09"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -6.91%  "

# Row 21
$ws.Range("E21").Value = "  -8.52%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "389.11"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.80%  "

# Row 23
$ws.Range("E23").Value = "  -8.39%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.35"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.43%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.17"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.31%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.91"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -6.26%  "

# Row 27
$ws.Range("E27").Value = "  -5.51%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.06"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.15%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.56"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -7.18%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.83"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -9.45%  "

# Row 31
$ws.Range("E31").Value = "  -3.85%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "30.83"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.49%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "629.46"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.50%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "12.16"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.00%  "

# Row 35
$ws.Range("E35").Value = "  -7.35%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "63.43"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.15%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "41.74"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -9.98%  "

# Row 38
$ws.Range("B38").Value = "Dai"
$ws.Range("C38").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.01"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.76%  "

# Row 39
$ws.Range("B39").Value = "TheGraph"
$ws.Range("C39").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.400"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.32%  "

# Row 40
$ws.Range("D40").Value = "0.0₃0749"
$ws.Range("E40").Value = "  -10.37%  "

# Row 41
$ws.Range("E41").Value = "  -5.57%  "

# Row 42
$ws.Range("D42").Value = "3.120.53"
$ws.Range("E42").Value = "  +6.39%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.07%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.93"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.69%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.64"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.19%  "

# Row 46
$ws.Range("E46").Value = "  -8.92%  "

# Row 47
$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.130"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -7.19%  "

# Row 48
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.08"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.50%  "

# Row 49
$ws.Range("E49").Value = "  -3.53%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.42"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -9.57%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.73"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.09%  "
